$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D's price strings to remain text (they can look numeric),
# matching the workbook's existing inline-string convention.
$priceRange = $ws.Range('D2:D51')
$priceRange.NumberFormat = '@'

$ws.Range('D2').Value = '22.452.63'
$ws.Range('E2').Value = '  +0.19%  '
$ws.Range('D3').Value = '1.573.87'
$ws.Range('E3').Value = '  +0.70%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('E5').Value = '  +0.04%  '
$ws.Range('D6').Value = '288.02'
$ws.Range('E6').Value = '  +0.78%  '
$ws.Range('D7').Value = '0.3723'
$ws.Range('E7').Value = '  +2.36%  '
$ws.Range('D8').Value = '47.56'
$ws.Range('E8').Value = '  -1.56%  '
$ws.Range('D9').Value = '0.3320'
$ws.Range('D10').Value = '1.156'
$ws.Range('E10').Value = '  +2.50%  '
$ws.Range('D11').Value = '0.07522'
$ws.Range('E11').Value = '  +1.46%  '
$ws.Range('E12').Value = '  +0.05%  '
$ws.Range('D13').Value = '20.79'
$ws.Range('E13').Value = '  -0.03%  '
$ws.Range('D14').Value = '5.934'
$ws.Range('E14').Value = '  +0.09%  '
$ws.Range('D15').Value = '6.928'
$ws.Range('E15').Value = '  +0.61%  '
$ws.Range('D16').Value = '1.567.23'
$ws.Range('E16').Value = '  +0.21%  '
$ws.Range('D17').Value = '0.00001117'
$ws.Range('E17').Value = '  +1.00%  '
$ws.Range('D18').Value = '88.34'
$ws.Range('E18').Value = '  -0.21%  '
$ws.Range('D19').Value = '0.06727'
$ws.Range('E19').Value = '  +0.46%  '
$ws.Range('E20').Value = '  +0.02%  '
$ws.Range('D21').Value = '6.393'
$ws.Range('E21').Value = '  +0.87%  '
$ws.Range('D22').Value = '16.51'
$ws.Range('E22').Value = '  +2.69%  '
$ws.Range('E23').Value = '  +0.33%  '
$ws.Range('D24').Value = '22.448.18'
$ws.Range('E24').Value = '  +0.19%  '
$ws.Range('D25').Value = '2.397'
$ws.Range('E25').Value = '  -0.77%  '
$ws.Range('D26').Value = '2.625'
$ws.Range('E26').Value = '  +3.00%  '
$ws.Range('D27').Value = '150.63'
$ws.Range('E27').Value = '  +0.62%  '
$ws.Range('D28').Value = '19.65'
$ws.Range('E28').Value = '  +1.39%  '
$ws.Range('D29').Value = '4.958'
$ws.Range('E29').Value = '  -0.83%  '
$ws.Range('D30').Value = '125.30'
$ws.Range('E30').Value = '  +1.87%  '
$ws.Range('D31').Value = '1.745.55'
$ws.Range('E31').Value = '  +0.47%  '
$ws.Range('D32').Value = '1.097'
$ws.Range('E32').Value = '  +3.01%  '
$ws.Range('D33').Value = '6.090'
$ws.Range('E33').Value = '  -0.47%  '
$ws.Range('D34').Value = '1.990'
$ws.Range('E34').Value = '  -0.40%  '
$ws.Range('D35').Value = '9.833'
$ws.Range('E35').Value = '  +2.41%  '
$ws.Range('D36').Value = '0.08337'
$ws.Range('E36').Value = '  +1.39%  '
$ws.Range('D37').Value = '0.02457'
$ws.Range('E37').Value = '  +2.67%  '
$ws.Range('D38').Value = '1.311'
$ws.Range('E38').Value = '  +0.67%  '
$ws.Range('D39').Value = '0.2233'
$ws.Range('E39').Value = '  +1.27%  '
$ws.Range('E40').Value = '  +0.15%  '
$ws.Range('E41').Value = '  +0.18%  '
$ws.Range('E42').Value = '  +2.51%  '
$ws.Range('D43').Value = '0.6269'
$ws.Range('E43').Value = '  +3.18%  '
$ws.Range('E44').Value = '  +0.05%  '
$ws.Range('D45').Value = '13.96'
$ws.Range('E45').Value = '  +2.31%  '
$ws.Range('D46').Value = '0.6091'
$ws.Range('E46').Value = '  +6.03%  '
$ws.Range('D48').Value = '2.047'
$ws.Range('E48').Value = '  +1.82%  '
$ws.Range('D49').Value = '125.02'
$ws.Range('E49').Value = '  +0.13%  '
$ws.Range('D50').Value = '1.210'
$ws.Range('E50').Value = '  -0.09%  '
$ws.Range('D51').Value = '0.07204'
$ws.Range('E51').Value = '  -0.13%  '

# Restore the default style so only the cell values changed (no formatting diff).
$priceRange.Style = 'Normal'

